$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new draw result as row 18, matching the existing rows which
# store every column as literal text (even the date-like and number-like
# values). A plain .Value assignment of "2025-10-04" / "251004" would be
# auto-coerced to a date serial / number by Excel's type inference, so we
# use a leading apostrophe to force literal text for those two, then copy
# the neighboring row's cell style back on top so no stray "quote prefix"
# number format is left behind on the new cells.

$ws.Cells.Item(18, 1).Value = "'2025-10-04"
$ws.Cells.Item(18, 1).Style = $ws.Cells.Item(17, 1).Style

$ws.Cells.Item(18, 2).Value = "Pick 3"

$ws.Cells.Item(18, 3).Value = "'251004"
$ws.Cells.Item(18, 3).Style = $ws.Cells.Item(17, 3).Style

$ws.Cells.Item(18, 4).Value = "0-6-8"

$ws.Cells.Item(18, 5).Value = "2025-10-04T21:34:57.554+04:00"
